# Updated attendance roll, created enemy class
# Fill in row 15 (Sprint 2, second meeting) of the attendance roll with the
# new meeting date/time, place, and attendance marks for each team member.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (font/border) of the row above onto C15 first, since
# filling the row visually "extends" the bordered/bold block of already
# completed rows (C6:C14) down by one row. This only changes the cell's
# style, not its value.
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New meeting date/time and place.
$ws.Range("B15").Value = "9/22 / 4:15"
$ws.Range("C15").Value = "Google Hangout"

# Attendance marks for the six team members (columns D-I).
$ws.Range("D15").Value = "A"
$ws.Range("E15").Value = "A"
$ws.Range("F15").Value = "A"
$ws.Range("G15").Value = "A"
$ws.Range("H15").Value = "A"
$ws.Range("I15").Value = "A"

# Update the selected/active cell as it was left after the edit.
$ws.Range("J15").Select() | Out-Null
